$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.388.20'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '3.425.91'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '406.38'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').Value = '130.14'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').Value = '  -2.63%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'0.690"
$ws.Range('E9').Value = '  +1.94%  '
$ws.Range('D10').Value = '0.137'
$ws.Range('E10').Value = '  +7.57%  '
$ws.Range('D11').Value = '41.93'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '19.81'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '8.37'
$ws.Range('E14').Value = '  -2.27%  '
$ws.Range('D15').Value = '3.397.19'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '62.300.61'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D18').Value = '1.01'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('E19').Value = '  +9.65%  '
$ws.Range('E20').Value = '  -3.07%  '
$ws.Range('D21').Value = '84.41'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').Value = '311.15'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '12.78'
$ws.Range('E23').Value = '  -3.07%  '
$ws.Range('D24').Value = '3.16'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').Value = '4.74'
$ws.Range('E25').Value = '  +3.63%  '
$ws.Range('D26').Value = '29.61'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = '8.07'
$ws.Range('E27').Value = '  -5.99%  '
$ws.Range('D28').Value = '7.74'
$ws.Range('E28').Value = '  +3.30%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.74'
$ws.Range('E29').Value = '  +4.54%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = "'44.90"
$ws.Range('E30').Value = '  +5.43%  '
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').Value = '51.84'
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').Value = '0.319'
$ws.Range('E39').Value = '  +11.07%  '
$ws.Range('D40').Value = "'3.30"
$ws.Range('D41').Value = '142.31'
$ws.Range('E41').Value = '  +3.79%  '
$ws.Range('E42').Value = '  -0.57%  '
$ws.Range('E43').Value = '  -3.65%  '
$ws.Range('D44').Value = "'3.90"
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('D45').Value = '16.79'
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('D47').Value = '21.08'
$ws.Range('E47').Value = '  -3.12%  '
$ws.Range('D48').Value = '2.102.10'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').Value = '1.98'
$ws.Range('E49').Value = '  +3.34%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '2.29'
$ws.Range('E50').Value = '  -2.25%  '
$ws.Range('B51').Value = 'OceanProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range('D51').Value = '1.09'
$ws.Range('E51').Value = '  +27.69%  '
